$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data blocks for the two report dates (rows 2-5 and rows 6-9) are
# swapped between each other, row by row, keeping the quality ("Calidad",
# column L) order fixed: Especial / Primera / Segunda / Tercera.
# Columns swapped per row pair: D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)

$pairs = @(
    @{A=2; B=6},
    @{A=3; B=7},
    @{A=4; B=8},
    @{A=5; B=9}
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $rowA = $pair.A
    $rowB = $pair.B

    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}
